# ---------------------------------------------------------------------------
# availabilities2.xlsx update
#   - tweak the "Conclusion" note
#   - re-left-align the whole sheet (was center-aligned)
#   - add a new "Parts / Cost / Link" cost breakdown table (rows 31-44)
#     covering a cinewhoop build, with a merged "Drone" label and two totals
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlLeft = -4131, xlCenter = -4108 (vertical)
$xlLeft = -4131
$xlCenterV = -4108

# --- 1. Small text tweak on the existing conclusion line ------------------
$ws.Range("B28").Value = "More costlier than I expected :(( around  1L"

# --- 2. Re-align everything currently on the sheet to left/center ---------
$ws.UsedRange.HorizontalAlignment = $xlLeft
$ws.UsedRange.VerticalAlignment = $xlCenterV

# --- 3. New cost-breakdown table ------------------------------------------

# Row 31 - note + link
$ws.Range("B31").Value = "Basic parts need to fly"
$ws.Range("D31").Value = "https://www.youtube.com/watch?v=QC-36Wfo0q4&list=PL_3XHkWVmS0L8qJ9UAu_fURkJFw8aD5k4&ab_channel=Techwittier-%E0%AE%A4%E0%AE%AE%E0%AE%BF%E0%AE%B4%E0%AF%8D"

# Row 32 - header row (bold)
$ws.Range("B32").Value = "Parts"
$ws.Range("C32").Value = "Cost"
$ws.Range("D32").Value = "Link"

# Row 33 - Radio
$ws.Range("A33").Value = "Radio"
$ws.Range("B33").Value = "Radio Transmitter with receiver and battery"
$ws.Range("C33").Value = 14000
$ws.Range("D33").Value = "https://robu.in/product/radiomaster-tx12-mkii-expresslrs-edgetx-radio-controller/?gad_source=1&gclid=Cj0KCQiAkKqsBhC3ARIsAEEjuJixrn-eq08sVZbY3zt8EQqFXp9a1I-jRgZi-mWJpKQS8F4uV-jFaN8aAozqEALw_wcB"

# Row 34 - Drone / Frame cinewhoop (first of the merged "Drone" block)
$ws.Range("A34").Value = "Drone"
$ws.Range("B34").Value = "Frame cinewhoop"
$ws.Range("C34").Value = 4000
$ws.Range("D34").Value = "https://www.quadkart.in/cloud-149-v2-3-inch-cinewhoop-frame-kit/"

# Row 35 - Motors
$ws.Range("B35").Value = "Motors (4) "
$ws.Range("C35").Value = 6000
$ws.Range("D35").Value = "https://www.drkstore.in/geprc-speedx2-1804-2450kv-3450kv-motor/"

# Row 36 - Stack FC and ESC
$ws.Range("B36").Value = "Stack FC and ESC"
$ws.Range("C36").Value = 6400
$ws.Range("D36").Value = "https://robu.in/product/f405-v1-0-flight-controller35a-2-6s-4-in-1-esc-flytower-mpu6000/?gad_source=1&gclid=Cj0KCQiAkKqsBhC3ARIsAEEjuJhtR27vIpTiLYbx6qre8xGnVikIKVSYDwGblV5JkCGVrV9qIF6jxEwaAsWsEALw_wcB"

# merge the "Drone" label across rows 34-36
$ws.Range("A34:A36").Merge()

# Row 37 - Batteries
$ws.Range("A37").Value = "Batteries"
$ws.Range("B37").Value = "1300mAh 4s 120C (4)"
$ws.Range("C37").Value = 10000
$ws.Range("D37").Value = "https://www.quadkart.in/tattu-r-line-version-3-0-1300mah-4s-120c-lipo-battery/"

# Row 38 - Charger
$ws.Range("A38").Value = "Charger"
$ws.Range("B38").Value = "SkyRC B6 Neo 200W - out of stock"
$ws.Range("C38").Value = 3500
$ws.Range("D38").Value = "https://www.drkstore.in/skyrc-b6-neo-200w-dc-smart-charger-with-dc-pd-dual-input/"

# Row 39 - Total (bold label)
$ws.Range("B39").Value = "Total"
$ws.Range("C39").Value = 45000

# Row 40 - follow-up note
$ws.Range("B40").Value = "To add FPV into above"

# Row 41 - VTX
$ws.Range("A41").Value = "VTX"
$ws.Range("B41").Value = "AKK Race Ranger VTX"
$ws.Range("C41").Value = 2200
$ws.Range("D41").Value = "https://www.quadkopters.com/product/fpv-kit-and-accessories/akk-race-ranger-vtx/"

# Row 42 - Camera
$ws.Range("A42").Value = "Camera"
$ws.Range("B42").Value = "1/3″ CMOS 1500TVL Mini FPV Camera 2.1mm Lens PAL / NTSC With OSD"
$ws.Range("C42").Value = 1500
$ws.Range("D42").Value = "https://robu.in/product/1-3-cmos-1500tvl-mini-fpv-camera-2-1mm-lens-pal-ntsc-with-osd/?gad_source=1&gclid=Cj0KCQiAkKqsBhC3ARIsAEEjuJh54Do8gbQ3zbHc6cWJ9_yUoc6Y-iCCvXei-JtfoILL0ASURxYIA4QaAnadEALw_wcB"

# Row 43 - Goggles
$ws.Range("A43").Value = "Goggles"
$ws.Range("B43").Value = "betafpv cetus"
$ws.Range("C43").Value = 15000

# Row 44 - grand Total (bold label)
$ws.Range("B44").Value = "Total "
$ws.Range("C44").Value = 65000

# --- 4. Bold the header / total rows, everything else stays regular -------
$boldRanges = "B20", "B27", "B32", "C32", "D32", "B39", "B44"
foreach ($r in $boldRanges) {
    $ws.Range($r).Font.Bold = $true
}

# Keep the new block left/center aligned as well (matches rest of sheet)
$ws.Range("A31:D44").HorizontalAlignment = $xlLeft
$ws.Range("A31:D44").VerticalAlignment = $xlCenterV

# --- 5. Selection / scroll position, like the saved file ------------------
$ws.Range("D43").Select()
